$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.712.19"
$ws.Range("E2").Value = "  -0.53%  "

$ws.Range("D3").Value = "2.275.33"
$ws.Range("E3").Value = "  -2.85%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.77%  "

$ws.Range("E7").Value = "  -1.08%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.510"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.53%  "

$ws.Range("E11").Value = "  -0.72%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.02%  "

$ws.Range("D14").Value = "2.622.57"
$ws.Range("E14").Value = "  -2.74%  "

$ws.Range("D15").Value = "2.271.43"
$ws.Range("E15").Value = "  -2.78%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.90%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.801"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.56%  "

$ws.Range("D18").Value = "46.618.89"
$ws.Range("E18").Value = "  -0.50%  "

$ws.Range("D19").Value = "0.0₃0994"
$ws.Range("E19").Value = "  +4.76%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.17%  "

$ws.Range("E21").Value = "  -5.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "247.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.72%  "

$ws.Range("E24").Value = "  -6.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.21%  "

$ws.Range("E27").Value = "  -1.74%  "

$ws.Range("E28").Value = "  -1.36%  "

$ws.Range("E29").Value = "  -2.71%  "

$ws.Range("E30").Value = "  -0.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +10.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "146.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.93%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.87%  "

$ws.Range("E35").Value = "  -5.00%  "

$ws.Range("E36").Value = "  +4.74%  "

$ws.Range("E37").Value = "  -2.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "15.69"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.17%  "

$ws.Range("E39").Value = "  -7.50%  "

$ws.Range("E40").Value = "  -4.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0298"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.55%  "

$ws.Range("E42").Value = "  -10.33%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "94.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +15.94%  "

$ws.Range("D45").Value = "1.787.09"
$ws.Range("E45").Value = "  -1.30%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "71.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.40%  "

$ws.Range("E48").Value = "  -6.49%  "

$ws.Range("E49").Value = "  -2.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "95.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.38%  "

$ws.Range("E51").Value = "  -0.93%  "
